$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Range("R33").Copy()
$ws.Range("Z4").Select()
$ws.Paste()
Write-Output "done"
